$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sD2 = $ws.Range("D2").Style
$ws.Range("D2").Value = "'61.225.03"
$ws.Range("D2").Style = $sD2
$ws.Range("E2").Value = "  +1.36%  "

$sD3 = $ws.Range("D3").Style
$ws.Range("D3").Value = "'3.347.12"
$ws.Range("D3").Style = $sD3
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  +0.10%  "

$sD5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'401.61"
$ws.Range("D5").Style = $sD5
$ws.Range("E5").Value = "  -1.78%  "

$sD6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'127.19"
$ws.Range("D6").Style = $sD6
$ws.Range("E6").Value = "  +11.88%  "

$sD7 = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.596"
$ws.Range("D7").Style = $sD7
$ws.Range("E7").Value = "  +4.30%  "

$sD8 = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = $sD8
$ws.Range("E8").Value = "  -0.06%  "

$sD9 = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.663"
$ws.Range("D9").Style = $sD9
$ws.Range("E9").Value = "  +6.37%  "

$sD10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.119"
$ws.Range("D10").Style = $sD10
$ws.Range("E10").Value = "  +3.87%  "

$sD11 = $ws.Range("D11").Style
$ws.Range("D11").Value = "'41.34"
$ws.Range("D11").Style = $sD11
$ws.Range("E11").Value = "  +5.01%  "

$ws.Range("E12").Value = "  -0.70%  "

$sD13 = $ws.Range("D13").Style
$ws.Range("D13").Value = "'3.901.10"
$ws.Range("D13").Style = $sD13
$ws.Range("E13").Value = "  +1.97%  "

$sD14 = $ws.Range("D14").Style
$ws.Range("D14").Value = "'8.37"
$ws.Range("D14").Style = $sD14
$ws.Range("E14").Value = "  +2.53%  "

$sD15 = $ws.Range("D15").Style
$ws.Range("D15").Value = "'19.40"
$ws.Range("D15").Style = $sD15
$ws.Range("E15").Value = "  +1.90%  "

$sD16 = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.393.17"
$ws.Range("D16").Style = $sD16
$ws.Range("E16").Value = "  +1.48%  "

$sD17 = $ws.Range("D17").Style
$ws.Range("D17").Value = "'61.133.72"
$ws.Range("D17").Style = $sD17
$ws.Range("E17").Value = "  +1.42%  "

$sD18 = $ws.Range("D18").Style
$ws.Range("D18").Value = "'11.26"
$ws.Range("D18").Style = $sD18
$ws.Range("E18").Value = "  +4.44%  "

$sD19 = $ws.Range("D19").Style
$ws.Range("D19").Value = "'1.01"
$ws.Range("D19").Style = $sD19
$ws.Range("E19").Value = "  +1.45%  "

$sD20 = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.0000129"
$ws.Range("D20").Style = $sD20
$ws.Range("E20").Value = "  +12.47%  "

$sD21 = $ws.Range("D21").Style
$ws.Range("D21").Value = "'3.21"
$ws.Range("D21").Style = $sD21
$ws.Range("E21").Value = "  -3.84%  "

$sD22 = $ws.Range("D22").Style
$ws.Range("D22").Value = "'81.16"
$ws.Range("D22").Style = $sD22
$ws.Range("E22").Value = "  +9.99%  "

$sD23 = $ws.Range("D23").Style
$ws.Range("D23").Value = "'12.92"
$ws.Range("D23").Style = $sD23
$ws.Range("E23").Value = "  +4.50%  "

$sD24 = $ws.Range("D24").Style
$ws.Range("D24").Value = "'303.03"
$ws.Range("D24").Style = $sD24
$ws.Range("E24").Value = "  +2.29%  "

$sD25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'3.20"
$ws.Range("D25").Style = $sD25
$ws.Range("E25").Value = "  +3.28%  "

$sD26 = $ws.Range("D26").Style
$ws.Range("D26").Value = "'4.68"
$ws.Range("D26").Style = $sD26
$ws.Range("E26").Value = "  +9.89%  "

$sD27 = $ws.Range("D27").Style
$ws.Range("D27").Value = "'8.40"
$ws.Range("D27").Style = $sD27
$ws.Range("E27").Value = "  +12.51%  "

$sD28 = $ws.Range("D28").Style
$ws.Range("D28").Value = "'29.19"
$ws.Range("D28").Style = $sD28
$ws.Range("E28").Value = "  +0.21%  "

$sD29 = $ws.Range("D29").Style
$ws.Range("D29").Value = "'7.42"
$ws.Range("D29").Style = $sD29
$ws.Range("E29").Value = "  -2.05%  "

$sD30 = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.172"
$ws.Range("D30").Style = $sD30
$ws.Range("E30").Value = "  +0.20%  "

$sD31 = $ws.Range("D31").Style
$ws.Range("D31").Value = "'0.114"
$ws.Range("D31").Style = $sD31
$ws.Range("E31").Value = "  +0.63%  "

$sD32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'11.51"
$ws.Range("D32").Style = $sD32
$ws.Range("E32").Value = "  +2.68%  "

$sD33 = $ws.Range("D33").Style
$ws.Range("D33").Value = "'2.55"
$ws.Range("D33").Style = $sD33
$ws.Range("E33").Value = "  +2.90%  "

$ws.Range("E34").Value = "  -0.07%  "

$sD35 = $ws.Range("D35").Style
$ws.Range("D35").Value = "'41.11"
$ws.Range("D35").Style = $sD35
$ws.Range("E35").Value = "  +1.78%  "

$sD36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.0479"
$ws.Range("D36").Style = $sD36
$ws.Range("E36").Value = "  -2.11%  "

$sD37 = $ws.Range("D37").Style
$ws.Range("D37").Value = "'52.19"
$ws.Range("D37").Style = $sD37
$ws.Range("E37").Value = "  +0.30%  "

$sD38 = $ws.Range("D38").Style
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = $sD38
$ws.Range("E38").Value = "  +0.22%  "

$sD39 = $ws.Range("D39").Style
$ws.Range("D39").Value = "'3.37"
$ws.Range("D39").Style = $sD39
$ws.Range("E39").Value = "  +0.91%  "

$sD40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.92"
$ws.Range("D40").Style = $sD40
$ws.Range("E40").Value = "  -4.11%  "

$sD41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'1.98"
$ws.Range("D41").Style = $sD41
$ws.Range("E41").Value = "  +5.44%  "

$sD42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.123"
$ws.Range("D42").Style = $sD42
$ws.Range("E42").Value = "  +2.97%  "

$sD43 = $ws.Range("D43").Style
$ws.Range("D43").Value = "'135.50"
$ws.Range("D43").Style = $sD43
$ws.Range("E43").Value = "  +0.92%  "

$sD44 = $ws.Range("D44").Style
$ws.Range("D44").Value = "'3.90"
$ws.Range("D44").Style = $sD44
$ws.Range("E44").Value = "  +3.22%  "

$sD45 = $ws.Range("D45").Style
$ws.Range("D45").Value = "'16.78"
$ws.Range("D45").Style = $sD45
$ws.Range("E45").Value = "  +3.38%  "

$sD46 = $ws.Range("D46").Style
$ws.Range("D46").Value = "'0.281"
$ws.Range("D46").Style = $sD46
$ws.Range("E46").Value = "  -3.93%  "

$sD47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.22"
$ws.Range("D47").Style = $sD47
$ws.Range("E47").Value = "  +1.59%  "

$sD48 = $ws.Range("D48").Style
$ws.Range("D48").Value = "'21.45"
$ws.Range("D48").Style = $sD48
$ws.Range("E48").Value = "  +2.91%  "

$sD49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'2.126.69"
$ws.Range("D49").Style = $sD49
$ws.Range("E49").Value = "  -0.25%  "

$sD50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'3.681.53"
$ws.Range("D50").Style = $sD50
$ws.Range("E50").Value = "  +1.66%  "

$sD51 = $ws.Range("D51").Style
$ws.Range("D51").Value = "'2.35"
$ws.Range("D51").Style = $sD51
$ws.Range("E51").Value = "  -1.09%  "
